$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "52.752.20"
$ws.Range("E2").Value2 = "  -13.01%  "
$ws.Range("D3").Value2 = "2.344.03"
$ws.Range("E3").Value2 = "  -19.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.00"
$ws.Range("E4").Value2 = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "447.55"
$ws.Range("E5").Value2 = "  -14.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "125.02"
$ws.Range("E6").Value2 = "  -12.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.999"
$ws.Range("E7").Value2 = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.473"
$ws.Range("E8").Value2 = "  -12.52%  "
$ws.Range("D9").Value2 = "2.347.47"
$ws.Range("E9").Value2 = "  -19.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.0912"
$ws.Range("E10").Value2 = "  -14.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "5.22"
$ws.Range("E11").Value2 = "  -14.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.309"
$ws.Range("E12").Value2 = "  -13.02%  "
$ws.Range("E13").Value2 = "  -6.07%  "
$ws.Range("D14").Value2 = "2.756.96"
$ws.Range("E14").Value2 = "  -19.00%  "
$ws.Range("D15").Value2 = "52.720.28"
$ws.Range("E15").Value2 = "  -13.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "19.13"
$ws.Range("E16").Value2 = "  -14.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "0.0000121"
$ws.Range("E17").Value2 = "  -13.53%  "
$ws.Range("D18").Value2 = "2.349.79"
$ws.Range("E18").Value2 = "  -19.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "4.09"
$ws.Range("E19").Value2 = "  -15.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "301.09"
$ws.Range("E20").Value2 = "  -13.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "9.02"
$ws.Range("E21").Value2 = "  -21.26%  "
$ws.Range("E22").Value2 = "  +1.01%  "
$ws.Range("E23").Value2 = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "5.26"
$ws.Range("E24").Value2 = "  -18.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "54.39"
$ws.Range("E25").Value2 = "  -15.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0.991"
$ws.Range("E26").Value2 = "  -0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "0.374"
$ws.Range("E27").Value2 = "  -16.21%  "
$ws.Range("D28").Value2 = "2.386.77"
$ws.Range("E28").Value2 = "  -21.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.142"
$ws.Range("E29").Value2 = "  -20.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "6.96"
$ws.Range("E30").Value2 = "  -10.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.998"
$ws.Range("E31").Value2 = "  -0.18%  "
$ws.Range("D32").Value2 = "0.0₃0677"
$ws.Range("E32").Value2 = "  -20.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "146.04"
$ws.Range("E33").Value2 = "  -4.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "17.25"
$ws.Range("E34").Value2 = "  -11.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.35"
$ws.Range("E35").Value2 = "  -18.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "4.86"
$ws.Range("E36").Value2 = "  -12.23%  "
$ws.Range("B37").Value2 = "Fetch.AI"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.808"
$ws.Range("E37").Value2 = "  -17.89%  "
$ws.Range("B38").Value2 = "NEARProtocol"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "3.40"
$ws.Range("E38").Value2 = "  -21.86%  "
$ws.Range("B39").Value2 = "ImmutableX"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "1.01"
$ws.Range("E39").Value2 = "  -14.53%  "
$ws.Range("B40").Value2 = "FirstDigitalUSD"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.995"
$ws.Range("E40").Value2 = "  -0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "32.65"
$ws.Range("E41").Value2 = "  -12.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.574"
$ws.Range("E42").Value2 = "  -11.72%  "
$ws.Range("B43").Value2 = "Filecoin"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "3.21"
$ws.Range("E43").Value2 = "  -12.54%  "
$ws.Range("B44").Value2 = "Hedera"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.0511"
$ws.Range("E44").Value2 = "  -11.57%  "
$ws.Range("B45").Value2 = "WhiteBITCoin"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "10.10"
$ws.Range("E45").Value2 = "  -2.31%  "
$ws.Range("D46").Value2 = "1.931.50"
$ws.Range("E46").Value2 = "  -15.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "1.18"
$ws.Range("E47").Value2 = "  -18.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.0209"
$ws.Range("E48").Value2 = "  -11.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.0836"
$ws.Range("E49").Value2 = "  -7.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "16.18"
$ws.Range("E50").Value2 = "  -19.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "4.09"
$ws.Range("E51").Value2 = "  -16.11%  "
